# Apply weekly fruit/vegetable data update: rows 2,3,4,5,9,10 get their
# D,K,L,M,N,O,P,Q,R,S,T values permuted (row 2<->9 swap; rows 3->4->5->10->3 cycle).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the update (D plus K..T)
$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

function Get-RowValues($ws, $row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value2()
    }
    return $vals
}

function Set-RowValues($ws, $row, $cols, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $vals[$c]
    }
}

# Snapshot current ("before") values for the affected rows
$row2  = Get-RowValues $ws 2  $cols
$row3  = Get-RowValues $ws 3  $cols
$row4  = Get-RowValues $ws 4  $cols
$row5  = Get-RowValues $ws 5  $cols
$row9  = Get-RowValues $ws 9  $cols
$row10 = Get-RowValues $ws 10 $cols

# Apply the permutation:
#   row 2  <- old row 9
#   row 9  <- old row 2
#   row 3  <- old row 4
#   row 4  <- old row 5
#   row 5  <- old row 10
#   row 10 <- old row 3
Set-RowValues $ws 2  $cols $row9
Set-RowValues $ws 9  $cols $row2
Set-RowValues $ws 3  $cols $row4
Set-RowValues $ws 4  $cols $row5
Set-RowValues $ws 5  $cols $row10
Set-RowValues $ws 10 $cols $row3
